# --------------------------------------------------------------------------
# Edit described by the commit:
#
#   1) The table on slide 16 (the "Google Shape;213;p29" graphicFrame) gets
#      a new table style applied:
#         {C9C74DD0-1D82-45F6-88AC-1A2C9086A982}  ->  {4B5BDF2F-E28E-4F64-B2F2-9E0E4E887C5E}
#
#   2) The deck's theme colours are swapped: the colour palette that used to
#      live in ppt/theme/theme1.xml ("Office Theme") becomes the palette
#      used by the presentation (ppt/theme/theme2.xml, the theme actually
#      referenced by the slide master / presentation.xml), i.e. the design
#      that was "Integral" now uses the stock Office colour values.
# --------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Table style on slide 16.
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(16)
$tableShape = $slide.Shapes.Item(3)
$table = $tableShape.Table
$table.ApplyStyle("{4B5BDF2F-E28E-4F64-B2F2-9E0E4E887C5E}")

# ---------------------------------------------------------------------
# 2) Theme colour scheme.
#
# The presentation's colour scheme (12 theme colours: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) is exposed through any slide's
# ThemeColorScheme and backs the single colour palette shared by every
# slide / the slide master. Push it to the "Office Theme" values.
# ---------------------------------------------------------------------
$themeColors = $slide.ThemeColorScheme

$themeColors.Item(1).RGB  = 0          # dk1      000000
$themeColors.Item(2).RGB  = 16777215   # lt1      FFFFFF
$themeColors.Item(3).RGB  = 6968388    # dk2      44546A
$themeColors.Item(4).RGB  = 15132391   # lt2      E7E6E6
$themeColors.Item(5).RGB  = 13998939   # accent1  5B9BD5
$themeColors.Item(6).RGB  = 3243501    # accent2  ED7D31
$themeColors.Item(7).RGB  = 10855845   # accent3  A5A5A5
$themeColors.Item(8).RGB  = 49407      # accent4  FFC000
$themeColors.Item(9).RGB  = 12874308   # accent5  4472C4
$themeColors.Item(10).RGB = 4697456    # accent6  70AD47
$themeColors.Item(11).RGB = 12673797   # hlink    0563C1
$themeColors.Item(12).RGB = 7491477    # folHlink 954F72
